$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 7166
$ws.Range("E2").Value = 528
$ws.Range("F2").Value = 528
$ws.Range("G2").Value = 496
$ws.Range("H2").Value = 380
$ws.Range("I2").Value = 406
$ws.Range("J2").Value = -25
$ws.Range("K2").Value = 7429
$ws.Range("L2").Value = 2698
$ws.Range("M2").Value = 4731
$ws.Range("N2").Value = 4711
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 409
$ws.Range("Q2").Value = 978
$ws.Range("R2").Value = -214
$ws.Range("S2").Value = -789
$ws.Range("T2").Value = 141
$ws.Range("U2").Value = 837
$ws.Range("V2").Value = 1830
$ws.Range("W2").Value = 7.36
$ws.Range("X2").Value = 5.31
$ws.Range("Y2").Value = 8.720000000000001
$ws.Range("Z2").Value = 5.01
$ws.Range("AA2").Value = 57.02
$ws.Range("AB2").Value = 1122.75
$ws.Range("AC2").Value = 993
$ws.Range("AD2").Value = 25.29
$ws.Range("AE2").Value = 11963
$ws.Range("AF2").Value = 2.1
$ws.Range("AG2").Value = 550
$ws.Range("AH2").Value = 2.19
$ws.Range("AI2").Value = 53.38
$ws.Range("AJ2").Value = 40878588
$ws.Range("D3").Value = 6004
$ws.Range("E3").Value = 436
$ws.Range("F3").Value = 436
$ws.Range("G3").Value = 421
$ws.Range("H3").Value = 342
$ws.Range("I3").Value = 366
$ws.Range("J3").Value = -24
$ws.Range("K3").Value = 7106
$ws.Range("L3").Value = 2377
$ws.Range("M3").Value = 4729
$ws.Range("N3").Value = 4730
$ws.Range("O3").Value = -1
$ws.Range("P3").Value = 409
$ws.Range("Q3").Value = 665
$ws.Range("R3").Value = -262
$ws.Range("S3").Value = -532
$ws.Range("T3").Value = 150
$ws.Range("U3").Value = 515
$ws.Range("V3").Value = 1678
$ws.Range("W3").Value = 7.26
$ws.Range("X3").Value = 5.7
$ws.Range("Y3").Value = 7.76
$ws.Range("Z3").Value = 4.71
$ws.Range("AA3").Value = 50.27
$ws.Range("AB3").Value = 1157.4
$ws.Range("AC3").Value = 896
$ws.Range("AD3").Value = 17.07
$ws.Range("AE3").Value = 12236
$ws.Range("AF3").Value = 1.25
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 3.27
$ws.Range("AI3").Value = 52.78
$ws.Range("AJ3").Value = 40878588
$ws.Range("D4").Value = 5982
$ws.Range("E4").Value = 769
$ws.Range("F4").Value = 769
$ws.Range("G4").Value = 591
$ws.Range("H4").Value = 399
$ws.Range("I4").Value = 431
$ws.Range("J4").Value = -32
$ws.Range("K4").Value = 7489
$ws.Range("L4").Value = 2556
$ws.Range("M4").Value = 4933
$ws.Range("N4").Value = 4931
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 409
$ws.Range("Q4").Value = 1155
$ws.Range("R4").Value = -580
$ws.Range("S4").Value = -510
$ws.Range("T4").Value = 176
$ws.Range("U4").Value = 979
$ws.Range("V4").Value = 1404
$ws.Range("W4").Value = 12.85
$ws.Range("X4").Value = 6.67
$ws.Range("Y4").Value = 8.93
$ws.Range("Z4").Value = 5.47
$ws.Range("AA4").Value = 51.81
$ws.Range("AB4").Value = 1213.7
$ws.Range("AC4").Value = 1055
$ws.Range("AD4").Value = 20.57
$ws.Range("AE4").Value = 12782
$ws.Range("AF4").Value = 1.7
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 2.3
$ws.Range("AI4").Value = 44.73
$ws.Range("AJ4").Value = 40878588
$ws.Range("D5").Value = 7158
$ws.Range("E5").Value = 1526
$ws.Range("F5").Value = 1526
$ws.Range("G5").Value = 1472
$ws.Range("H5").Value = 1101
$ws.Range("I5").Value = 1144
$ws.Range("J5").Value = -43
$ws.Range("K5").Value = 8816
$ws.Range("L5").Value = 2828
$ws.Range("M5").Value = 5988
$ws.Range("N5").Value = 5809
$ws.Range("O5").Value = 179
$ws.Range("P5").Value = 409
$ws.Range("Q5").Value = 1564
$ws.Range("R5").Value = -1715
$ws.Range("S5").Value = 155
$ws.Range("T5").Value = 195
$ws.Range("U5").Value = 1369
$ws.Range("V5").Value = 1420
$ws.Range("W5").Value = 21.32
$ws.Range("X5").Value = 15.38
$ws.Range("Y5").Value = 21.3
$ws.Range("Z5").Value = 13.5
$ws.Range("AA5").Value = 47.22
$ws.Range("AB5").Value = 1453.02
$ws.Range("AC5").Value = 2798
$ws.Range("AD5").Value = 8.77
$ws.Range("AE5").Value = 15161
$ws.Range("AF5").Value = 1.62
$ws.Range("AG5").Value = 1500
$ws.Range("AH5").Value = 6.11
$ws.Range("AI5").Value = 50.24
$ws.Range("AJ5").Value = 40878588
$ws.Range("D6").Value = 7629
$ws.Range("E6").Value = 1392
$ws.Range("F6").Value = 1392
$ws.Range("G6").Value = 1407
$ws.Range("H6").Value = 1051
$ws.Range("I6").Value = 1045
$ws.Range("K6").Value = 8647
$ws.Range("L6").Value = 1987
$ws.Range("M6").Value = 6661
$ws.Range("N6").Value = 6474
$ws.Range("P6").Value = 409
$ws.Range("Q6").Value = 1166
$ws.Range("R6").Value = -191
$ws.Range("S6").Value = -748
$ws.Range("T6").Value = 242
$ws.Range("U6").Value = 924
$ws.Range("V6").Value = 941
$ws.Range("W6").Value = 18.25
$ws.Range("X6").Value = 13.77
$ws.Range("Y6").Value = 17.02
$ws.Range("Z6").Value = 12.04
$ws.Range("AA6").Value = 29.83
$ws.Range("AB6").Value = 1580.28
$ws.Range("AC6").Value = 2556
$ws.Range("AD6").Value = 9.43
$ws.Range("AE6").Value = 16488
$ws.Range("AF6").Value = 1.46
$ws.Range("AG6").Value = 1200
$ws.Range("AH6").Value = 4.98
$ws.Range("AI6").Value = 45.09
$ws.Range("AJ6").Value = 40878588
$ws.Range("D7").Value = 6570
$ws.Range("E7").Value = 1155
$ws.Range("G7").Value = 1158
$ws.Range("H7").Value = 836
$ws.Range("I7").Value = 834
$ws.Range("K7").Value = 9181
$ws.Range("L7").Value = 2192
$ws.Range("M7").Value = 6990
$ws.Range("N7").Value = 6780
$ws.Range("P7").Value = 409
$ws.Range("Q7").Value = 1112
$ws.Range("R7").Value = -338
$ws.Range("S7").Value = -426
$ws.Range("T7").Value = 268
$ws.Range("U7").Value = 870
$ws.Range("W7").Value = 17.58
$ws.Range("X7").Value = 12.72
$ws.Range("Y7").Value = 12.58
$ws.Range("Z7").Value = 9.380000000000001
$ws.Range("AA7").Value = 31.36
$ws.Range("AC7").Value = 2040
$ws.Range("AD7").Value = 9.56
$ws.Range("AE7").Value = 17670
$ws.Range("AF7").Value = 1.1
$ws.Range("AG7").Value = 1033
$ws.Range("AH7").Value = 5.3
$ws.Range("AI7").Value = 50.65
$ws.Range("D8").Value = 6664
$ws.Range("E8").Value = 1058
$ws.Range("G8").Value = 1052
$ws.Range("H8").Value = 788
$ws.Range("I8").Value = 787
$ws.Range("K8").Value = 9476
$ws.Range("L8").Value = 2114
$ws.Range("M8").Value = 7362
$ws.Range("N8").Value = 7135
$ws.Range("P8").Value = 409
$ws.Range("Q8").Value = 1036
$ws.Range("R8").Value = -646
$ws.Range("S8").Value = -417
$ws.Range("T8").Value = 341
$ws.Range("U8").Value = 534
$ws.Range("W8").Value = 15.88
$ws.Range("X8").Value = 11.82
$ws.Range("Y8").Value = 11.31
$ws.Range("Z8").Value = 8.44
$ws.Range("AA8").Value = 28.72
$ws.Range("AC8").Value = 1925
$ws.Range("AD8").Value = 10.13
$ws.Range("AE8").Value = 18595
$ws.Range("AF8").Value = 1.05
$ws.Range("AG8").Value = 1033
$ws.Range("AH8").Value = 5.3
$ws.Range("AI8").Value = 53.69
$ws.Range("D9").Value = 6573
$ws.Range("E9").Value = 1157
$ws.Range("G9").Value = 1159
$ws.Range("H9").Value = 868
$ws.Range("I9").Value = 865
$ws.Range("K9").Value = 9887
$ws.Range("L9").Value = 2136
$ws.Range("M9").Value = 7750
$ws.Range("N9").Value = 7486
$ws.Range("P9").Value = 409
$ws.Range("Q9").Value = 1137
$ws.Range("R9").Value = -804
$ws.Range("S9").Value = -341
$ws.Range("T9").Value = 367
$ws.Range("U9").Value = 533
$ws.Range("W9").Value = 17.6
$ws.Range("X9").Value = 13.21
$ws.Range("Y9").Value = 11.83
$ws.Range("Z9").Value = 8.970000000000001
$ws.Range("AA9").Value = 27.56
$ws.Range("AC9").Value = 2116
$ws.Range("AD9").Value = 9.210000000000001
$ws.Range("AE9").Value = 19508
$ws.Range("AF9").Value = 1
$ws.Range("AG9").Value = 1050
$ws.Range("AH9").Value = 5.38
$ws.Range("AI9").Value = 49.61
